$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - card holder first name
$ws.Range("C2").Value = "Hartmut"

# Row 3 - card number (must stay text) and surname
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Row 5 - opening balance date
$ws.Range("D5").Value = "KONTOSTAND AM 24.05.2024"

# Row 6 - transaction 1
$ws.Range("B6").Value = "25.05."
$ws.Range("C6").Value = "26.05."
$ws.Range("D6").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E6").Value = "24,84-"

# Row 7 - transaction 2
$ws.Range("B7").Value = "27.05."
$ws.Range("C7").Value = "28.05."
$ws.Range("D7").Value = "BEITRAG Allianz SE K-65873143"
$ws.Range("E7").Value = "57,68-"

# Row 8 - transaction 3
$ws.Range("B8").Value = "30.05."
$ws.Range("C8").Value = "31.05."
$ws.Range("D8").Value = "ZALANDO MKTPLC EU GMKPEH"
$ws.Range("E8").Value = "62,00-"

# Row 12 - closing balance date and amount
$ws.Range("D12").Value = "KONTOSTAND AM 04.06.2024"
$ws.Range("E12").Value = "144,52-"

# Row 13 - next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 11.06.2024"
